$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("C1").Value = "inssurance"

# Update age (column B) and eligibility/insurance (column C) values for rows 2-20
$ages = @(21, 45, 32, 41, 20, 35, 20, 23, 42, 34, 24, 22, 23, 25, 43, 44, 25, 30, 31)
$labels = @("no", "yes", "yes", "yes", "no", "yes", "no", "no", "yes", "yes", "no", "no", "no", "no", "yes", "yes", "no", "yes", "yes")

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $ages[$i]
    $ws.Cells.Item($row, 3).Value = $labels[$i]
}

# Remove the old row 21 (was A21=19, B21=44, C21=0), data now ends at row 20
$ws.Rows.Item(21).Delete()
